$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coinranking snapshot refresh: update the Price (D) and Volume(1h) (E) columns
# for every listed coin row (2-51) with the latest scraped figures.
#
# A handful of new Price values are plain decimal numbers (e.g. '313.72',
# '0.579', '15.70'). Left alone, Excel's smart-entry parsing would coerce
# those into numeric cells (and even drop significant trailing zeros, e.g.
# '15.70' -> 15.7), whereas the sheet stores Price/Volume as plain text.
# Forcing NumberFormat="@" (Text) before the assignment keeps them text,
# matching every other cell in these columns (multi-dot prices like
# '42.598.61' or percent strings like '  -0.58%  ' are already safe as-is).

$ws.Range("D2").Value = '42.598.61'
$ws.Range("E2").Value = '  -0.58%  '
$ws.Range("D3").Value = '2.539.05'
$ws.Range("E3").Value = '  -0.17%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '313.72'
$ws.Range("E5").Value = '  +3.35%  '
$ws.Range("E6").Value = '  -2.68%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.579'
$ws.Range("E7").Value = '  +0.45%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.538'
$ws.Range("E9").Value = '  -1.52%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.17'
$ws.Range("E10").Value = '  -1.70%  '
$ws.Range("E11").Value = '  -1.54%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.69'
$ws.Range("E12").Value = '  -0.39%  '
$ws.Range("E13").Value = '  -0.48%  '
$ws.Range("D14").Value = '2.933.58'
$ws.Range("E14").Value = '  -0.15%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.70'
$ws.Range("E15").Value = '  +4.61%  '
$ws.Range("D16").Value = '2.533.49'
$ws.Range("E16").Value = '  -1.22%  '
$ws.Range("E17").Value = '  -0.58%  '
$ws.Range("D18").Value = '42.660.45'
$ws.Range("E18").Value = '  -0.58%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.08'
$ws.Range("E19").Value = '  -1.44%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.67'
$ws.Range("E20").Value = '  +1.54%  '
$ws.Range("D21").Value = '0.0₃0969'
$ws.Range("E21").Value = '  -2.01%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '71.08'
$ws.Range("E22").Value = '  -1.15%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '254.84'
$ws.Range("E23").Value = '  -0.23%  '
$ws.Range("E24").Value = '  +0.17%  '
$ws.Range("E25").Value = '  -1.40%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '27.67'
$ws.Range("E26").Value = '  -1.45%  '
$ws.Range("E27").Value = '  +0.16%  '
$ws.Range("E28").Value = '  +1.91%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '39.69'
$ws.Range("E29").Value = '  +4.94%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '10.06'
$ws.Range("E30").Value = '  -1.03%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.97'
$ws.Range("E31").Value = '  -1.80%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '155.34'
$ws.Range("E32").Value = '  -1.68%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.41'
$ws.Range("E33").Value = '  +2.89%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '19.52'
$ws.Range("E34").Value = '  +0.47%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.14'
$ws.Range("E35").Value = '  +0.65%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0791'
$ws.Range("E36").Value = '  -0.81%  '
$ws.Range("E37").Value = '  -0.48%  '
$ws.Range("E38").Value = '  -4.08%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '24.81'
$ws.Range("E39").Value = '  -3.90%  '
$ws.Range("E40").Value = '  -0.25%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.24'
$ws.Range("E41").Value = '  +7.14%  '
$ws.Range("E42").Value = '  -0.91%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.84'
$ws.Range("E43").Value = '  -0.98%  '
$ws.Range("E44").Value = '  -0.74%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.999'
$ws.Range("E45").Value = '  -0.07%  '
$ws.Range("D46").Value = '2.050.94'
$ws.Range("E46").Value = '  -1.89%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '85.08'
$ws.Range("E47").Value = '  -3.22%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.92'
$ws.Range("E48").Value = '  -0.03%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '75.37'
$ws.Range("E49").Value = '  +0.63%  '
$ws.Range("D50").Value = '2.788.05'
$ws.Range("E50").Value = '  -0.31%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.191'
$ws.Range("E51").Value = '  +0.01%  '
